$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (header stays the same in row 1)
$data = @(
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Taurean Prince", "SG,SF", "Milwaukee Bucks"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jalen Johnson", "PF", "Atlanta Hawks"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
